$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row of data: Date (A27) and Error Count (B27)
$ws.Range("A26").Copy()
$ws.Range("A27").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A27").Value = 45995
$ws.Range("B27").Value = 69

# Move the active selection, matching the author's final cursor position
$ws.Range("E23").Select()
